$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing Q column values (detect_structure resets to 0) ---
$ws.Range("Q55").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("Q63").Value = 0
$ws.Range("Q71").Value = 0
$ws.Range("Q74").Value = 0
$ws.Range("Q79").Value = 0

# --- Fix R366 / O367 / R367 (inlineStr blank -> numeric 0 / update isPivot) ---
$ws.Range("R366").Value = 0
$ws.Range("O367").Value = 1
$ws.Range("R367").Value = 0

# --- Append new weekly rows 368-397 ---
# Row 368
$ws.Range("A368").Value = 45474
$ws.Range("A368").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B368").Value = 1208.491887850962
$ws.Range("C368").Value = 1267.742897965314
$ws.Range("D368").Value = 1143.073525326799
$ws.Range("E368").Value = 1160.02734375
$ws.Range("G368").Value = 35947290
$ws.Range("H368").Value = 2024
$ws.Range("I368").Value = 7
$ws.Range("J368").Value = 1
$ws.Range("K368").Value = 0
$ws.Range("L368").Value = 0
$ws.Range("M368").Value = 0
$ws.Range("N368").Value = 27
$ws.Range("O368").Value = 0
$ws.Range("P368").Value = 0
$ws.Range("Q368").Value = 0

# Row 369
$ws.Range("A369").Value = 45481
$ws.Range("A369").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B369").Value = 1168.541723734418
$ws.Range("C369").Value = 1245.945128511698
$ws.Range("D369").Value = 1128.591579333241
$ws.Range("E369").Value = 1207.942504882812
$ws.Range("G369").Value = 20972746
$ws.Range("H369").Value = 2024
$ws.Range("I369").Value = 7
$ws.Range("J369").Value = 8
$ws.Range("K369").Value = 0
$ws.Range("L369").Value = 0
$ws.Range("M369").Value = 0
$ws.Range("N369").Value = 28
$ws.Range("O369").Value = 0
$ws.Range("P369").Value = 0
$ws.Range("Q369").Value = 0

# Row 370
$ws.Range("A370").Value = 45488
$ws.Range("A370").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B370").Value = 1213.98506924921
$ws.Range("C370").Value = 1235.807812256028
$ws.Range("D370").Value = 1136.331972470826
$ws.Range("E370").Value = 1141.175903320312
$ws.Range("G370").Value = 7016704
$ws.Range("H370").Value = 2024
$ws.Range("I370").Value = 7
$ws.Range("J370").Value = 15
$ws.Range("K370").Value = 0
$ws.Range("L370").Value = 0
$ws.Range("M370").Value = 0
$ws.Range("N370").Value = 29
$ws.Range("O370").Value = 0
$ws.Range("P370").Value = 0
$ws.Range("Q370").Value = 0

# Row 371
$ws.Range("A371").Value = 45495
$ws.Range("A371").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B371").Value = 1132.625
$ws.Range("C371").Value = 1234.75
$ws.Range("D371").Value = 1085.5
$ws.Range("E371").Value = 1223.400024414062
$ws.Range("G371").Value = 11194778
$ws.Range("H371").Value = 2024
$ws.Range("I371").Value = 7
$ws.Range("J371").Value = 22
$ws.Range("K371").Value = 0
$ws.Range("L371").Value = 0
$ws.Range("M371").Value = 0
$ws.Range("N371").Value = 30
$ws.Range("O371").Value = 2
$ws.Range("P371").Value = 0
$ws.Range("Q371").Value = 0

# Row 372
$ws.Range("A372").Value = 45502
$ws.Range("A372").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B372").Value = 1235.625
$ws.Range("C372").Value = 1297.099975585938
$ws.Range("D372").Value = 1195.650024414062
$ws.Range("E372").Value = 1226.599975585938
$ws.Range("G372").Value = 16893192
$ws.Range("H372").Value = 2024
$ws.Range("I372").Value = 7
$ws.Range("J372").Value = 29
$ws.Range("K372").Value = 0
$ws.Range("L372").Value = 0
$ws.Range("M372").Value = 0
$ws.Range("N372").Value = 31
$ws.Range("O372").Value = 0
$ws.Range("P372").Value = 0
$ws.Range("Q372").Value = 0

# Row 373
$ws.Range("A373").Value = 45509
$ws.Range("A373").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B373").Value = 1225
$ws.Range("C373").Value = 1297.425048828125
$ws.Range("D373").Value = 1139.599975585938
$ws.Range("E373").Value = 1282.625
$ws.Range("G373").Value = 24270024
$ws.Range("H373").Value = 2024
$ws.Range("I373").Value = 8
$ws.Range("J373").Value = 5
$ws.Range("K373").Value = 0
$ws.Range("L373").Value = 0
$ws.Range("M373").Value = 0
$ws.Range("N373").Value = 32
$ws.Range("O373").Value = 0
$ws.Range("P373").Value = 0
$ws.Range("Q373").Value = 2

# Row 374
$ws.Range("A374").Value = 45516
$ws.Range("A374").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B374").Value = 1295.449951171875
$ws.Range("C374").Value = 1406
$ws.Range("D374").Value = 1254
$ws.Range("E374").Value = 1394.074951171875
$ws.Range("G374").Value = 22389648
$ws.Range("H374").Value = 2024
$ws.Range("I374").Value = 8
$ws.Range("J374").Value = 12
$ws.Range("K374").Value = 0
$ws.Range("L374").Value = 0
$ws.Range("M374").Value = 0
$ws.Range("N374").Value = 33
$ws.Range("O374").Value = 0
$ws.Range("P374").Value = 0
$ws.Range("Q374").Value = 0

# Row 375
$ws.Range("A375").Value = 45523
$ws.Range("A375").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B375").Value = 1427.5
$ws.Range("C375").Value = 1664.400024414062
$ws.Range("D375").Value = 1420
$ws.Range("E375").Value = 1568.5
$ws.Range("G375").Value = 43023968
$ws.Range("H375").Value = 2024
$ws.Range("I375").Value = 8
$ws.Range("J375").Value = 19
$ws.Range("K375").Value = 0
$ws.Range("L375").Value = 0
$ws.Range("M375").Value = 0
$ws.Range("N375").Value = 34
$ws.Range("O375").Value = 1
$ws.Range("P375").Value = 0
$ws.Range("Q375").Value = 0

# Row 376
$ws.Range("A376").Value = 45530
$ws.Range("A376").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B376").Value = 1600
$ws.Range("C376").Value = 1615
$ws.Range("D376").Value = 1405
$ws.Range("E376").Value = 1434.900024414062
$ws.Range("G376").Value = 19665903
$ws.Range("H376").Value = 2024
$ws.Range("I376").Value = 8
$ws.Range("J376").Value = 26
$ws.Range("K376").Value = 0
$ws.Range("L376").Value = 0
$ws.Range("M376").Value = 0
$ws.Range("N376").Value = 35
$ws.Range("O376").Value = 0
$ws.Range("P376").Value = 0
$ws.Range("Q376").Value = 0

# Row 377
$ws.Range("A377").Value = 45537
$ws.Range("A377").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B377").Value = 1443
$ws.Range("C377").Value = 1454
$ws.Range("D377").Value = 1400
$ws.Range("E377").Value = 1401.5
$ws.Range("G377").Value = 8046916
$ws.Range("H377").Value = 2024
$ws.Range("I377").Value = 9
$ws.Range("J377").Value = 2
$ws.Range("K377").Value = 0
$ws.Range("L377").Value = 0
$ws.Range("M377").Value = 0
$ws.Range("N377").Value = 36
$ws.Range("O377").Value = 0
$ws.Range("P377").Value = 0
$ws.Range("Q377").Value = 1

# Row 378
$ws.Range("A378").Value = 45544
$ws.Range("A378").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B378").Value = 1400
$ws.Range("C378").Value = 1410.75
$ws.Range("D378").Value = 1348.400024414062
$ws.Range("E378").Value = 1381.449951171875
$ws.Range("G378").Value = 8704760
$ws.Range("H378").Value = 2024
$ws.Range("I378").Value = 9
$ws.Range("J378").Value = 9
$ws.Range("K378").Value = 0
$ws.Range("L378").Value = 0
$ws.Range("M378").Value = 0
$ws.Range("N378").Value = 37
$ws.Range("O378").Value = 0
$ws.Range("P378").Value = 0
$ws.Range("Q378").Value = 0

# Row 379
$ws.Range("A379").Value = 45551
$ws.Range("A379").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B379").Value = 1392.5
$ws.Range("C379").Value = 1581
$ws.Range("D379").Value = 1372.099975585938
$ws.Range("E379").Value = 1538.650024414062
$ws.Range("G379").Value = 35008101
$ws.Range("H379").Value = 2024
$ws.Range("I379").Value = 9
$ws.Range("J379").Value = 16
$ws.Range("K379").Value = 0
$ws.Range("L379").Value = 0
$ws.Range("M379").Value = 0
$ws.Range("N379").Value = 38
$ws.Range("O379").Value = 0
$ws.Range("P379").Value = 0
$ws.Range("Q379").Value = 0

# Row 380
$ws.Range("A380").Value = 45558
$ws.Range("A380").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B380").Value = 1554
$ws.Range("C380").Value = 1559.800048828125
$ws.Range("D380").Value = 1454.800048828125
$ws.Range("E380").Value = 1463.050048828125
$ws.Range("G380").Value = 12259003
$ws.Range("H380").Value = 2024
$ws.Range("I380").Value = 9
$ws.Range("J380").Value = 23
$ws.Range("K380").Value = 0
$ws.Range("L380").Value = 0
$ws.Range("M380").Value = 0
$ws.Range("N380").Value = 39
$ws.Range("O380").Value = 0
$ws.Range("P380").Value = 0
$ws.Range("Q380").Value = 0

# Row 381
$ws.Range("A381").Value = 45565
$ws.Range("A381").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B381").Value = 1462.400024414062
$ws.Range("C381").Value = 1482.550048828125
$ws.Range("D381").Value = 1366.199951171875
$ws.Range("E381").Value = 1373.800048828125
$ws.Range("G381").Value = 10186949
$ws.Range("H381").Value = 2024
$ws.Range("I381").Value = 9
$ws.Range("J381").Value = 30
$ws.Range("K381").Value = 0
$ws.Range("L381").Value = 0
$ws.Range("M381").Value = 0
$ws.Range("N381").Value = 40
$ws.Range("O381").Value = 0
$ws.Range("P381").Value = 0
$ws.Range("Q381").Value = 2

# Row 382
$ws.Range("A382").Value = 45572
$ws.Range("A382").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B382").Value = 1385
$ws.Range("C382").Value = 1515
$ws.Range("D382").Value = 1344.599975585938
$ws.Range("E382").Value = 1478.050048828125
$ws.Range("G382").Value = 27245474
$ws.Range("H382").Value = 2024
$ws.Range("I382").Value = 10
$ws.Range("J382").Value = 7
$ws.Range("K382").Value = 0
$ws.Range("L382").Value = 0
$ws.Range("M382").Value = 0
$ws.Range("N382").Value = 41
$ws.Range("O382").Value = 2
$ws.Range("P382").Value = 0
$ws.Range("Q382").Value = 0

# Row 383
$ws.Range("A383").Value = 45579
$ws.Range("A383").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B383").Value = 1485
$ws.Range("C383").Value = 1640
$ws.Range("D383").Value = 1480.349975585938
$ws.Range("E383").Value = 1579.300048828125
$ws.Range("G383").Value = 46300166
$ws.Range("H383").Value = 2024
$ws.Range("I383").Value = 10
$ws.Range("J383").Value = 14
$ws.Range("K383").Value = 0
$ws.Range("L383").Value = 0
$ws.Range("M383").Value = 0
$ws.Range("N383").Value = 42
$ws.Range("O383").Value = 0
$ws.Range("P383").Value = 0
$ws.Range("Q383").Value = 0

# Row 384
$ws.Range("A384").Value = 45586
$ws.Range("A384").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B384").Value = 1592.599975585938
$ws.Range("C384").Value = 1604.900024414062
$ws.Range("D384").Value = 1381.199951171875
$ws.Range("E384").Value = 1419
$ws.Range("G384").Value = 17714360
$ws.Range("H384").Value = 2024
$ws.Range("I384").Value = 10
$ws.Range("J384").Value = 21
$ws.Range("K384").Value = 0
$ws.Range("L384").Value = 0
$ws.Range("M384").Value = 0
$ws.Range("N384").Value = 43
$ws.Range("O384").Value = 0
$ws.Range("P384").Value = 0
$ws.Range("Q384").Value = 0

# Row 385
$ws.Range("A385").Value = 45593
$ws.Range("A385").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B385").Value = 1472
$ws.Range("C385").Value = 1570
$ws.Range("D385").Value = 1447.699951171875
$ws.Range("E385").Value = 1558
$ws.Range("G385").Value = 20227314
$ws.Range("H385").Value = 2024
$ws.Range("I385").Value = 10
$ws.Range("J385").Value = 28
$ws.Range("K385").Value = 0
$ws.Range("L385").Value = 0
$ws.Range("M385").Value = 0
$ws.Range("N385").Value = 44
$ws.Range("O385").Value = 0
$ws.Range("P385").Value = 0
$ws.Range("Q385").Value = 0

# Row 386
$ws.Range("A386").Value = 45600
$ws.Range("A386").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B386").Value = 1556.449951171875
$ws.Range("C386").Value = 1678.849975585938
$ws.Range("D386").Value = 1504
$ws.Range("E386").Value = 1596.75
$ws.Range("G386").Value = 25765498
$ws.Range("H386").Value = 2024
$ws.Range("I386").Value = 11
$ws.Range("J386").Value = 4
$ws.Range("K386").Value = 0
$ws.Range("L386").Value = 0
$ws.Range("M386").Value = 0
$ws.Range("N386").Value = 45
$ws.Range("O386").Value = 0
$ws.Range("P386").Value = 0
$ws.Range("Q386").Value = 0

# Row 387
$ws.Range("A387").Value = 45607
$ws.Range("A387").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B387").Value = 1587.900024414062
$ws.Range("C387").Value = 1589.75
$ws.Range("D387").Value = 1427
$ws.Range("E387").Value = 1472.75
$ws.Range("G387").Value = 14347327
$ws.Range("H387").Value = 2024
$ws.Range("I387").Value = 11
$ws.Range("J387").Value = 11
$ws.Range("K387").Value = 0
$ws.Range("L387").Value = 0
$ws.Range("M387").Value = 0
$ws.Range("N387").Value = 46
$ws.Range("O387").Value = 0
$ws.Range("P387").Value = 0
$ws.Range("Q387").Value = 0

# Row 388
$ws.Range("A388").Value = 45614
$ws.Range("A388").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B388").Value = 1480
$ws.Range("C388").Value = 1553.699951171875
$ws.Range("D388").Value = 1451.449951171875
$ws.Range("E388").Value = 1540.25
$ws.Range("G388").Value = 12410411
$ws.Range("H388").Value = 2024
$ws.Range("I388").Value = 11
$ws.Range("J388").Value = 18
$ws.Range("K388").Value = 0
$ws.Range("L388").Value = 0
$ws.Range("M388").Value = 0
$ws.Range("N388").Value = 47
$ws.Range("O388").Value = 0
$ws.Range("P388").Value = 0
$ws.Range("Q388").Value = 1

# Row 389
$ws.Range("A389").Value = 45621
$ws.Range("A389").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B389").Value = 1560
$ws.Range("C389").Value = 1663.900024414062
$ws.Range("D389").Value = 1550
$ws.Range("E389").Value = 1639.449951171875
$ws.Range("G389").Value = 18805451
$ws.Range("H389").Value = 2024
$ws.Range("I389").Value = 11
$ws.Range("J389").Value = 25
$ws.Range("K389").Value = 0
$ws.Range("L389").Value = 0
$ws.Range("M389").Value = 0
$ws.Range("N389").Value = 48
$ws.Range("O389").Value = 0
$ws.Range("P389").Value = 0
$ws.Range("Q389").Value = 0

# Row 390
$ws.Range("A390").Value = 45628
$ws.Range("A390").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B390").Value = 1647
$ws.Range("C390").Value = 1899
$ws.Range("D390").Value = 1630
$ws.Range("E390").Value = 1883.800048828125
$ws.Range("G390").Value = 41981248
$ws.Range("H390").Value = 2024
$ws.Range("I390").Value = 12
$ws.Range("J390").Value = 2
$ws.Range("K390").Value = 0
$ws.Range("L390").Value = 0
$ws.Range("M390").Value = 0
$ws.Range("N390").Value = 49
$ws.Range("O390").Value = 0
$ws.Range("P390").Value = 0
$ws.Range("Q390").Value = 0

# Row 391
$ws.Range("A391").Value = 45635
$ws.Range("A391").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B391").Value = 1875
$ws.Range("C391").Value = 1962.199951171875
$ws.Range("D391").Value = 1859.699951171875
$ws.Range("E391").Value = 1943.5
$ws.Range("G391").Value = 18747542
$ws.Range("H391").Value = 2024
$ws.Range("I391").Value = 12
$ws.Range("J391").Value = 9
$ws.Range("K391").Value = 0
$ws.Range("L391").Value = 0
$ws.Range("M391").Value = 0
$ws.Range("N391").Value = 50
$ws.Range("O391").Value = 0
$ws.Range("P391").Value = 0
$ws.Range("Q391").Value = 0

# Row 392
$ws.Range("A392").Value = 45642
$ws.Range("A392").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B392").Value = 1948.599975585938
$ws.Range("C392").Value = 1989.800048828125
$ws.Range("D392").Value = 1851.599975585938
$ws.Range("E392").Value = 1860.75
$ws.Range("G392").Value = 14618649
$ws.Range("H392").Value = 2024
$ws.Range("I392").Value = 12
$ws.Range("J392").Value = 16
$ws.Range("K392").Value = 0
$ws.Range("L392").Value = 0
$ws.Range("M392").Value = 0
$ws.Range("N392").Value = 51
$ws.Range("O392").Value = 1
$ws.Range("P392").Value = 0
$ws.Range("Q392").Value = 0

# Row 393
$ws.Range("A393").Value = 45649
$ws.Range("A393").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B393").Value = 1878
$ws.Range("C393").Value = 1888
$ws.Range("D393").Value = 1762.550048828125
$ws.Range("E393").Value = 1777.550048828125
$ws.Range("G393").Value = 8892186
$ws.Range("H393").Value = 2024
$ws.Range("I393").Value = 12
$ws.Range("J393").Value = 23
$ws.Range("K393").Value = 0
$ws.Range("L393").Value = 0
$ws.Range("M393").Value = 0
$ws.Range("N393").Value = 52
$ws.Range("O393").Value = 0
$ws.Range("P393").Value = 0
$ws.Range("Q393").Value = 0

# Row 394
$ws.Range("A394").Value = 45656
$ws.Range("A394").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B394").Value = 1779
$ws.Range("C394").Value = 1828.5
$ws.Range("D394").Value = 1720.349975585938
$ws.Range("E394").Value = 1805.349975585938
$ws.Range("G394").Value = 10826148
$ws.Range("H394").Value = 2024
$ws.Range("I394").Value = 12
$ws.Range("J394").Value = 30
$ws.Range("K394").Value = 0
$ws.Range("L394").Value = 0
$ws.Range("M394").Value = 0
$ws.Range("N394").Value = 1
$ws.Range("O394").Value = 0
$ws.Range("P394").Value = 0
$ws.Range("Q394").Value = 0

# Row 395
$ws.Range("A395").Value = 45663
$ws.Range("A395").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B395").Value = 1807.400024414062
$ws.Range("C395").Value = 1816.900024414062
$ws.Range("D395").Value = 1590
$ws.Range("E395").Value = 1594.449951171875
$ws.Range("G395").Value = 12585025
$ws.Range("H395").Value = 2025
$ws.Range("I395").Value = 1
$ws.Range("J395").Value = 6
$ws.Range("K395").Value = 0
$ws.Range("L395").Value = 0
$ws.Range("M395").Value = 0
$ws.Range("N395").Value = 2
$ws.Range("O395").Value = 0
$ws.Range("P395").Value = 0
$ws.Range("Q395").Value = 0

# Row 396
$ws.Range("A396").Value = 45670
$ws.Range("A396").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B396").Value = 1574.849975585938
$ws.Range("C396").Value = 1642.25
$ws.Range("D396").Value = 1516.75
$ws.Range("E396").Value = 1597.900024414062
$ws.Range("G396").Value = 14572323
$ws.Range("H396").Value = 2025
$ws.Range("I396").Value = 1
$ws.Range("J396").Value = 13
$ws.Range("K396").Value = 0
$ws.Range("L396").Value = 0
$ws.Range("M396").Value = 0
$ws.Range("N396").Value = 3
$ws.Range("O396").Value = 0
$ws.Range("P396").Value = 0
$ws.Range("Q396").Value = 0

# Row 397
$ws.Range("A397").Value = 45677
$ws.Range("A397").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B397").Value = 1607.699951171875
$ws.Range("C397").Value = 1619.349975585938
$ws.Range("D397").Value = 1474.050048828125
$ws.Range("E397").Value = 1500.25
$ws.Range("G397").Value = 12325670
$ws.Range("H397").Value = 2025
$ws.Range("I397").Value = 1
$ws.Range("J397").Value = 20
$ws.Range("K397").Value = 0
$ws.Range("L397").Value = 0
$ws.Range("M397").Value = 0
$ws.Range("N397").Value = 4
$ws.Range("O397").Value = 0
$ws.Range("P397").Value = 0
$ws.Range("Q397").Value = 0

Write-Output "done"